$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells sometimes hold numeric-looking text such as
# "1.00" or "155.00". If we just assign a plain string, Excel's normal
# type inference turns them into real numbers and the literal text is
# lost. To keep them as text (matching the source data), force each
# target cell to Text format before writing the value, then restore the
# default (General) cell style afterwards so no formatting changes leak
# into the file.
$dRows = @(2,3,5,6,7,8,9,11,14,15,16,18,19,20,21,22,23,24,25,26,27,28,29,30,32,33,35,36,38,39,40,41,42,43,44,45,46,47,48,50,51)
foreach ($r in $dRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '60.452.02'
$ws.Range("E2").Value = '  -3.07%  '
$ws.Range("D3").Value = '2.580.79'
$ws.Range("E3").Value = '  -4.45%  '
$ws.Range("E4").Value = '  +0.59%  '
$ws.Range("D5").Value = '508.37'
$ws.Range("E5").Value = '  -2.32%  '
$ws.Range("D6").Value = '155.00'
$ws.Range("E6").Value = '  -4.46%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("D8").Value = '0.580'
$ws.Range("E8").Value = '  -5.98%  '
$ws.Range("D9").Value = '2.588.56'
$ws.Range("E9").Value = '  -4.24%  '
$ws.Range("E10").Value = '  +6.40%  '
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  -3.22%  '
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").Value = '3.037.70'
$ws.Range("E14").Value = '  -2.26%  '
$ws.Range("D15").Value = '60.478.21'
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '21.65'
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("D18").Value = '2.594.66'
$ws.Range("E18").Value = '  -3.08%  '
$ws.Range("D19").Value = '4.76'
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("D20").Value = '345.71'
$ws.Range("E20").Value = '  -5.66%  '
$ws.Range("D21").Value = '10.47'
$ws.Range("E21").Value = '  -2.22%  '
$ws.Range("D22").Value = '6.11'
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").Value = '60.11'
$ws.Range("E24").Value = '  -1.63%  '
$ws.Range("D25").Value = '0.420'
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").Value = '2.704.09'
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +1.62%  '
$ws.Range("D29").Value = '0.0₃0845'
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("D30").Value = '7.38'
$ws.Range("E30").Value = '  -3.74%  '
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").Value = '19.38'
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("D33").Value = '152.54'
$ws.Range("E33").Value = '  -4.08%  '
$ws.Range("E34").Value = '  -2.62%  '
$ws.Range("D35").Value = '5.71'
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").Value = '4.01'
$ws.Range("E36").Value = '  -1.19%  '
$ws.Range("E37").Value = '  -3.68%  '
$ws.Range("D38").Value = '0.850'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("D39").Value = '1.48'
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("D40").Value = '0.846'
$ws.Range("E40").Value = '  -3.85%  '
$ws.Range("D41").Value = '36.13'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").Value = '3.75'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").Value = '296.69'
$ws.Range("D44").Value = '0.623'
$ws.Range("E44").Value = '  -4.75%  '
$ws.Range("D45").Value = '0.0997'
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("D46").Value = '0.0558'
$ws.Range("E46").Value = '  -5.13%  '
$ws.Range("D47").Value = '0.998'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '19.78'
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("E49").Value = '  -4.35%  '
$ws.Range("D50").Value = '0.0233'
$ws.Range("E50").Value = '  -3.21%  '
$ws.Range("D51").Value = '10.29'
$ws.Range("E51").Value = '  +0.23%  '

foreach ($r in $dRows) {
    $ws.Range("D$r").Style = "Normal"
}

